$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an exact text string into a cell without altering its
# style/number-format (direct .Value assignment on a numeric-looking
# string would silently convert it to a Double and drop formatting like
# trailing zeros, e.g. "24.70" -> 24.7). Routing the text through a
# formula cell (="...") and pasting values-only preserves the literal
# text while leaving the destination cells style untouched.
$helperCell = $ws.Range("H1")
function Set-ExactText($cellRef, $text) {
    $helperCell.Formula = '="' + $text + '"'
    $helperCell.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helperCell.Clear()
}

# --- Price column (D): set as exact text to match source formatting ---
Set-ExactText "D2" '30.268.61'
Set-ExactText "D3" '2.087.73'
Set-ExactText "D5" '343.28'
Set-ExactText "D7" '0.5215'
Set-ExactText "D8" '0.4405'
Set-ExactText "D9" '54.48'
Set-ExactText "D10" '0.09328'
Set-ExactText "D12" '24.70'
Set-ExactText "D13" '8.653'
Set-ExactText "D14" '6.890'
Set-ExactText "D15" '2.065.24'
Set-ExactText "D16" '100.99'
Set-ExactText "D17" '0.00001154'
Set-ExactText "D18" '1.006'
Set-ExactText "D20" '0.06681'
Set-ExactText "D21" '6.361'
Set-ExactText "D22" '1.004'
Set-ExactText "D23" '30.238.28'
Set-ExactText "D24" '12.49'
Set-ExactText "D26" '21.73'
Set-ExactText "D27" '162.22'
Set-ExactText "D28" '2.507'
Set-ExactText "D29" '132.69'
Set-ExactText "D32" '0.1046'
Set-ExactText "D33" '6.215'
Set-ExactText "D34" '6.646'
Set-ExactText "D36" '10.11'
Set-ExactText "D37" '0.02622'
Set-ExactText "D38" '0.06774'
Set-ExactText "D39" '0.6971'
Set-ExactText "D40" '1.342'
Set-ExactText "D41" '12.50'
Set-ExactText "D42" '0.2212'
Set-ExactText "D43" '0.6806'
Set-ExactText "D44" '14.28'
Set-ExactText "D45" '2.330'
Set-ExactText "D48" '3.630'
Set-ExactText "D49" '0.00000000346'
Set-ExactText "D50" '1.212'

# --- Volume(1h) column (E): plain text assignment (never numeric-looking) ---
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  +5.37%  '
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("E21").Value = '  +2.77%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("E34").Value = '  +10.31%  '
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  +4.24%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("E44").Value = '  +0.73%  '
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("E47").Value = '  +18.30%  '
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("E50").Value = '  +8.43%  '
$ws.Range("E51").Value = '  -0.31%  '
